$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "New York" value got attached one row too early (to row 34's "New York"
# label) while the true "New York" value had slid down into a spurious
# "New York City" row. Delete that extra "New York City" row entirely -
# this removes the duplicate label, shifts every following state up by one
# row, and re-aligns each label with its correct labor-force value.
$ws.Rows.Item(35).Delete()

# Leave the view/selection where it landed after the delete (whole row 35,
# which now holds the state that slid up into that slot).
$ws.Rows.Item(35).Select()
